# CategoryTab.xlsx update ([02/16] Excel2Json 시스템 제작)
#
# - Rename the worksheet from "시트1" to "Data"
# - The "sword" category row (id 10001001) was missing its tab_name value;
#   fill it in so the row matches the pattern of the other categories
#   (shield, bow, armor, food, etc).
# - Leave the cursor/selection on G10.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Name = "Data"

$ws.Range("B3").Value = "sword"

$ws.Range("G10").Select()
